$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / coin-name updates (never numeric-looking, safe to assign directly)
$ws.Range("D2").Value = '66.049.34'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '2.690.46'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("E6").Value = '  +1.94%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E10").Value = '  +4.34%  '
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("E13").Value = '  +9.60%  '
$ws.Range("E14").Value = '  +3.12%  '
$ws.Range("D15").Value = '3.176.20'
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("D16").Value = '65.924.02'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").Value = '2.691.82'
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("E18").Value = '  +1.42%  '
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("E20").Value = '  +6.44%  '
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("E22").Value = '  +3.60%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  +18.10%  '
$ws.Range("E25").Value = '  +5.27%  '
$ws.Range("E26").Value = '  -1.45%  '
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("E28").Value = '  +3.92%  '
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("E30").Value = '  +2.36%  '
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("E32").Value = '  +2.79%  '
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +1.97%  '
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  +1.54%  '
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("E45").Value = '  +2.13%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E47").Value = '  +4.57%  '
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("E50").Value = '  +6.43%  '
$ws.Range("E51").Value = '  +1.34%  '

# Numeric-looking text values: force text storage via quote-prefix, then restore the
# default "Normal" cell style so no stray number-format/style gets attached.
$ws.Range("D5").Value = "'612.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'158.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Value = "'6.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.404"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Value = "'30.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Value = "'12.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'4.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'7.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'358.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'71.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'9.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'1.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'0.172"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Value = "'2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Value = "'540.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Value = "'6.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'5.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.436"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'20.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'164.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Value = "'168.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'42.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'4.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.0632"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'23.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Value = "'0.659"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'20.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0995"
$ws.Range("D51").Style = "Normal"
